$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1100, 10.23025608062744, 1.984598636627197, 0.8545932173728943, 19.19128799438477, 18476, "02"),
    @(1200, 10.20132637023926, 2.017646074295044, 0.8375071883201599, 18.7402172088623,  18568, "02"),
    @(1300, 10.15693473815918, 2.02303671836853,  0.8238384127616882, 18.05677795410156, 18438, "02"),
    @(1400, 10.11833477020264, 2.039085626602173, 1.051082253456116,  17.29816055297852, 18366, "02"),
    @(1500, 10.08100700378418, 2.092822313308716, 1.206564784049988,  19.58426666259766, 18392, "02")
)

$startRow = 21

# Ensure the Month column (G) is formatted as text so values like "02" keep
# their leading zero instead of being coerced to the number 2.
$ws.Range("G$startRow`:G25").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Output "Added rows 21-25"
